$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1029804860, "Info Pembelian token", "2023-08-27 12:24:08.809666"),
    @(3, 1029804860, "Info Pembelian token", "2023-08-27 12:29:22.081420"),
    @(4, 1029804860, "Info Pembelian token", "2023-08-28 11:56:19.713636"),
    @(5, 1029804860, "Info Pembelian token", "2023-08-28 15:24:38.833037"),
    @(6, 1029804860, "Info Pembelian token", "2023-08-28 15:25:10.724502"),
    @(7, 1029804860, "Info Pembelian token", "2023-08-28 15:27:37.693104"),
    @(8, 1029804860, "Info Pembelian token", "2023-08-28 15:30:30.944016"),
    @(9, 1029804860, "Info Pembelian token", "2023-08-28 15:35:04.575610"),
    @(10, 1029804860, "Info Pelanggan", "2023-08-28 15:35:47.626011"),
    @(11, 1029804860, "Info Pembelian token", "2023-08-28 15:36:20.304000"),
    @(12, 1029804860, "Info Pembelian token", "2023-08-28 15:36:41.239902"),
    @(13, 1029804860, "Info Pembelian token", "2023-08-28 16:08:24.880039"),
    @(14, 1029804860, "Info Pembelian token", "2023-08-28 16:19:35.497789"),
    @(15, 1029804860, "Info Pembelian token", "2023-08-28 20:27:41.937119"),
    @(16, 1029804860, "Info Pembelian token", "2023-08-28 20:31:20.912933"),
    @(17, 1029804860, "Info Pembelian token", "2023-08-28 20:32:58.325586"),
    @(18, 1029804860, "Info Pembelian token", "2023-08-28 20:51:46.004077"),
    @(19, 1029804860, "Info Pembelian token", "2023-08-28 21:11:03.332767")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
